$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new blank column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before N (existing N..P shift right to O..Q)
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of the column to its left (M)
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Select/activate cell R7 and make this sheet the active one (matches the
# recorded selection/active-tab state in the workbook)
$ws.Range("R7").Select()
$ws.Activate()

$wb.Save()
